# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.611.39'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.583.43'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.47'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.36'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.572.14'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.620'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.215'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +17.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.651'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.26'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000324'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +6.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.56'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.148.95'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.60'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.539.08'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.573.51'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.42'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '567.27'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +15.60%  '
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.02'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.74'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -6.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.67'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.49%  '
$ws.Range('E25').Value = '  -0.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '95.51'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.52'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.15'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.31'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.32'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.48'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.88'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.32'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '563.99'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.419'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.91'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.21%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0787'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.375.51'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.29%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.134'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.06'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.55'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0446'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.97'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.50'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.84%  '
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.997'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.46'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -10.37%  '

Write-Output "Applied 96 cell updates"
